$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F, matching style of existing header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Per-row "time_taken" timestamps for rows 2-29
$timeTaken = @(
    "2021-10-05 10:51:22.448073",
    "2021-10-05 10:51:22.448083",
    "2021-10-05 10:51:22.448087",
    "2021-10-05 10:51:22.448090",
    "2021-10-05 10:51:22.448093",
    "2021-10-05 10:51:22.448096",
    "2021-10-05 10:51:22.448098",
    "2021-10-05 10:51:22.448101",
    "2021-10-05 10:51:22.448104",
    "2021-10-05 10:51:22.448107",
    "2021-10-05 10:51:22.448110",
    "2021-10-05 10:51:22.448112",
    "2021-10-05 10:51:22.448115",
    "2021-10-05 10:51:22.448118",
    "2021-10-05 10:51:22.448121",
    "2021-10-05 10:51:22.448123",
    "2021-10-05 10:51:22.448126",
    "2021-10-05 10:51:22.448129",
    "2021-10-05 10:51:22.448132",
    "2021-10-05 10:51:22.448135",
    "2021-10-05 10:51:22.448137",
    "2021-10-05 10:51:22.448140",
    "2021-10-05 10:51:22.448143",
    "2021-10-05 10:51:22.448146",
    "2021-10-05 10:51:22.448148",
    "2021-10-05 10:51:22.448151",
    "2021-10-05 10:51:22.448154",
    "2021-10-05 10:51:22.448157"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
